$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 ("I0") and J1 ("IF"), matching H1's header style ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-36: new I column values, J column is a copy of H ---
$iValues = 5,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,4,1
$jValues = 7,6,4,6,5,7,7,6,9,6,6,5,5,3,7,6,6,5,7,6,6,7,7,7,7,5,5,7,6,6,6,5,4,6,2

for ($r = 2; $r -le 36; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
